$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Original column layout: A=label, B=value, C=minimum, D=maximum, E=vary,
# F=non-negative, G=expression
#
# Target column layout:   A=label, B=value, C=expression, D=minimum,
# E=maximum, F=non-negative, G=vary, H=standard-error
#
# Insert a fresh column before C. This shifts minimum/maximum/vary/
# non-negative/expression one column to the right (and keeps the empty
# "minimum"/"maximum" cells intact), landing them at D/E/F/G/H.
$ws.Columns("C:C").Insert()

# After the insert the layout is:
# A=label, B=value, C=(new, empty), D=minimum, E=maximum, F=vary,
# G=non-negative, H=expression

# Find the last data row (header is row 1).
$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row

# Fix up the header row.
$ws.Cells.Item(1, 3).Value = "expression"
$ws.Cells.Item(1, 6).Value = "non-negative"
$ws.Cells.Item(1, 7).Value = "vary"
$ws.Cells.Item(1, 8).Value = "standard-error"

for ($r = 2; $r -le $lastRow; $r++) {
    # Capture the values that are about to be overwritten.
    $varyValue = $ws.Cells.Item($r, 6).Value2
    $nonNegativeValue = $ws.Cells.Item($r, 7).Value2
    $expressionValue = $ws.Cells.Item($r, 8).Value2

    # C: expression (moved from the old "expression" column).
    $ws.Cells.Item($r, 3).Value = $expressionValue

    # F/G: non-negative/vary swap places.
    $ws.Cells.Item($r, 6).Value = $nonNegativeValue
    $ws.Cells.Item($r, 7).Value = $varyValue

    # H: new standard-error column, always "None" in this dataset.
    $ws.Cells.Item($r, 8).Value = "None"
}
